$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1. Wording fix inside the "Custom library" explanation bullet: a comma is
#    added after "library".
# -------------------------------------------------------------------------
$d.Content.Find.Execute(
    "has the library it needs and then installs them",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "has the library, it needs and then installs them", 2) | Out-Null

# -------------------------------------------------------------------------
# 2. New content at the end of the list: an "items" bullet with a hitbox
#    write-up, and an "Enemies" bullet that repeats the custom-library
#    write-up (original wording, no comma).
# -------------------------------------------------------------------------
$pTop    = $d.Paragraphs(20)   # "Custom library " - top-level bullet (numId 14)
$pDetail = $d.Paragraphs(21)   # "...library, it needs...[fixed by jaden] " (numId 15)

# (a) bare blank paragraph right after the custom-library explanation
$pDetail.Range.InsertParagraphAfter()
$pBlank1 = $d.Paragraphs($pDetail.Index + 1)
$pBlank1.Range.ListFormat.RemoveNumbers()
$pBlank1.TabStops.ClearAll()
$pBlank1.Style = "Normal"

# (b) "items " top-level bullet (duplicate pTop's numbering format)
$pTop.Range.InsertParagraphAfter()
$pItems = $d.Paragraphs($pTop.Index + 1)
$pItems.Range.Text = "items "

# (c) hitbox detail bullet (duplicate pDetail's numbering format)
$pDetail.Range.InsertParagraphAfter()
$pHitbox = $d.Paragraphs($pDetail.Index + 1)
$pHitbox.Range.Text = "I had problems with the hitboxes as I was making my own hitbox system, but it was very Janky, I realised that I could just apply a rectangle to the tombs with a list holding there coordinates and use the default hitbox system [fixed by jaden]"

# (d) indented blank paragraph after the hitbox bullet
$pHitbox.Range.InsertParagraphAfter()
$pBlank2 = $d.Paragraphs($pHitbox.Index + 1)
$pBlank2.Range.ListFormat.RemoveNumbers()
$pBlank2.TabStops.ClearAll()
$pBlank2.Style = "Normal"
$pBlank2.LeftIndent = $word.InchesToPoints(0.5)

# (e) "Enemies " top-level bullet
$pTop.Range.InsertParagraphAfter()
$pEnemies = $d.Paragraphs($pTop.Index + 1)
$pEnemies.Range.Text = "Enemies "

# (f) duplicated custom-library explanation bullet (original wording)
$pDetail.Range.InsertParagraphAfter()
$pEnemiesDetail = $d.Paragraphs($pDetail.Index + 1)
$pEnemiesDetail.Range.Text = "I had a problem with other people running the code as they don" + [char]0x2019 + "t have the liberty" + [char]0x2019 + "s, I now do a check the first time it runs to see if it has the library it needs and then installs them if they don" + [char]0x2019 + "t [fixed by jaden] "

# (g) indented blank paragraph after the duplicated bullet
$pEnemiesDetail.Range.InsertParagraphAfter()
$pBlank3 = $d.Paragraphs($pEnemiesDetail.Index + 1)
$pBlank3.Range.ListFormat.RemoveNumbers()
$pBlank3.TabStops.ClearAll()
$pBlank3.Style = "Normal"
$pBlank3.LeftIndent = $word.InchesToPoints(0.5)

Write-Output "edit complete"
